# Split Vensim views into submodules (#283)
#
# On the "No monotonous" worksheet a second block of sample data (rows 10-14,
# mirroring the existing rows 3-7 block) is appended, a new merged cell
# (A12:A14) is created, the sheet's used range grows to A1:J14, and the
# worksheet becomes the active sheet/tab of the workbook with C11 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("No monotonous")

# Duplicate the existing data block (rows 3:7, including values, number
# types, shared-string cells, cell styles and the A5:A7 merged cell) down
# into the new rows 10:14.
$ws.Range("A3:J7").Copy()
$ws.Range("A10:J14").PasteSpecial(-4104)

# A handful of values differ between the original block and the new one.
$ws.Range("G10").Value = 4
$ws.Range("I11").Value = 4
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 2

# Rows 6/7 (and now 13/14) only contain an empty, merged, styled A cell.
# Touching the alignment materializes that otherwise-empty formatted cell
# so it is written out (and keeps the sheet's dimension/used range in
# sync), mirroring A6/A7.
$ws.Range("A13").HorizontalAlignment = -4108
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("A14").VerticalAlignment = -4108

# Make "No monotonous" the active sheet/tab with C11 selected.
$ws.Activate()
$ws.Range("C11").Select()
